# Update the "2020" Advent of Code personal stats sheet with y20d4 results.
# Day 3 and Day 4 rows previously only had placeholder labels ("Day 3: " /
# "Day 4: ") with no timing data; now they're filled in with the real
# puzzle titles and completion times (matching the already-populated rows
# above/below them).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("2020")

# Row 7 -> Day 3
$ws.Range("B7").Value = "Day 3: Toboggan Trajectory"
$ws.Range("C7").Value = 0.3833333333333333
$ws.Range("E7").Value = 0.52361111111111114
$ws.Range("F7").Value = 0.26319444444444445
$ws.Range("H7").Value = "7th"

# Row 8 -> Day 4
$ws.Range("B8").Value = "Day 4: Passport Processing"
$ws.Range("C8").Value = 0.54861111111111105
$ws.Range("E8").Value = 1.3041666666666667
$ws.Range("F8").Value = 0.76666666666666661
$ws.Range("H8").Value = "8th"

# Move the visible selection to H9, matching the committed workbook state.
$ws.Range("H9").Select()
